# Questao 3 - marks the correct checkbox-style answer ("style.css, apenas.")
# by turning the blank red run inside "[      ]" into a red "X", and wraps
# the inserted "X" with grammar-check proofing markers, exactly like the
# already-answered paragraphs elsewhere in the document.

$d = $word.ActiveDocument

# Locate the target list paragraph: "[      ] style.css, apenas."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("[      ] style.css, apenas.")) {
        $target = $para
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the '[      ] style.css, apenas.' answer paragraph"
}

# Replace the whole paragraph (pPr + runs, including the trailing paragraph
# mark) via WordOpenXML so the four existing runs ("[  ", the two FF0000
# blank-space runs, "]") are preserved untouched except that the first
# FF0000 run's text becomes "X" (instead of two spaces) and is bracketed by
# <w:proofErr gramStart/.../gramEnd/> markers - replacing a sub-range in
# place reliably only works when the InsertXML target spans a full
# paragraph (Start..End), otherwise the engine re-appends content at the
# paragraph's end instead of in place.
$full = $d.Range($target.Range.Start, $target.Range.End)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="65B541B1" w14:textId="2D111261" w:rsidR="00FE1DFC" w:rsidRDefault="005B540F" w:rsidP="00951497"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">[  </w:t></w:r><w:r w:rsidR="00A83CDC"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>X</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00A53210"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>]</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FE1DFC"><w:t>style.css, apenas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$full.InsertXML($xml)
